# HungryDragonContent_General.xlsx
# "deleted tkhwbk movie seasonal" -- remove the "movie_kwwbk" season row from
# the {seasonsDefinitions} table on the global_settings sheet, and re-focus
# the workbook on that sheet / that row (matching the author's final
# selection state after making the edit in the Excel UI).

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("global_settings")
$wsTech = $wb.Worksheets.Item("tech")

# The "movie_kwwbk" entry lives on row 21 of the seasonsDefinitions table
# (sku/icon/tidName = movie_kwwbk / icon_season_movie_kwwbk /
# TID_SEASON_MOVIE_KWWBK_NAME). Delete the whole row; everything below
# (remaining seasons, pre-registration rewards table, etc.) shifts up by one
# row automatically, along with the table/autofilter ranges, conditional
# formatting and data validation that reference them.
$wsSettings.Rows.Item(21).Delete()

# Re-create the author's final view state: global_settings becomes the
# active/selected sheet (it was "tech" before), with the newly-shifted row 21
# selected in full.
$wsSettings.Activate()
$wsSettings.Range("A21:XFD21").Select()
